$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the other header cells (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the I0 / IF data columns for rows 2-42
$data = @(
    @(2,7,7),
    @(3,7,8),
    @(4,9,9),
    @(5,6,7),
    @(6,9,9),
    @(7,7,7),
    @(8,7,7),
    @(9,7,7),
    @(10,8,8),
    @(11,6,6),
    @(12,6,6),
    @(13,5,5),
    @(14,8,9),
    @(15,7,8),
    @(16,4,4),
    @(17,6,6),
    @(18,6,6),
    @(19,9,9),
    @(20,5,5),
    @(21,5,6),
    @(22,8,8),
    @(23,9,9),
    @(24,5,5),
    @(25,7,7),
    @(26,7,7),
    @(27,8,8),
    @(28,8,8),
    @(29,1,2),
    @(30,6,7),
    @(31,5,5),
    @(32,8,8),
    @(33,9,9),
    @(34,7,8),
    @(35,6,6),
    @(36,7,7),
    @(37,5,5),
    @(38,4,4),
    @(39,6,6),
    @(40,5,6),
    @(41,7,7),
    @(42,4,4)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $iF = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $iF
}
